# Remove the stray "Ι.Κ.Υ." bullet paragraph that sits right after the
# "ΚΟΙΝΟΠΟΙΗΣΗ" heading (just before the "${local_directorate}" bullet).
# The paragraph's own text is split across two runs ("Ι" + ".Κ.Υ. "),
# so locate it with Find and delete the whole paragraph (including its
# end-of-paragraph mark) while leaving the following paragraph untouched.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Ι.Κ.Υ.", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Expand the hit to its enclosing paragraph (wdParagraph = 4) so the
    # paragraph mark is included and the whole paragraph disappears.
    $rng.Expand(4) | Out-Null
    $rng.Delete()
}
